$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap columns B and C: cut C, insert before B
$ws.Columns("C:C").Cut() | Out-Null
$ws.Columns("B:B").Insert() | Out-Null

# Swap columns D and E: cut E, insert before D
$ws.Columns("E:E").Cut() | Out-Null
$ws.Columns("D:D").Insert() | Out-Null

# Select column D as the final selection (entire column)
$ws.Columns("D:D").Select() | Out-Null
